# Append the new match row (row 63) to the Azerbaijan Premier League sheet,
# mirroring the formatting of the previous data row (row 62) and filling in
# the new match's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 62's formatting (styles/number formats) down into row 63 first,
# so the new row picks up the same cellXfs (bold/bordered index column,
# date-time number format on the match-date column) as every other row.
$ws.Range("A62:V62").Copy()
$ws.Range("A63:V63").PasteSpecial(-4122)

# Now populate the new row's values.
$ws.Range("A63").Value = 62
$ws.Range("B63").Value = "azerbaijan"
$ws.Range("C63").Value = "premier-league"
$ws.Range("D63").Value = "2023-2024"
$ws.Range("E63").Value = 45241.66666666666
$ws.Range("F63").Value = "Turan"
$ws.Range("G63").Value = 2
$ws.Range("H63").Value = "Sabah Baku"
$ws.Range("I63").Value = 3
$ws.Range("J63").Value = 2.95
$ws.Range("K63").Value = "10/11/2023 04:12"
$ws.Range("L63").Value = 3.62
$ws.Range("M63").Value = "11/11/2023 15:19"
$ws.Range("N63").Value = 3.52
$ws.Range("O63").Value = "10/11/2023 04:12"
$ws.Range("P63").Value = 3.56
$ws.Range("Q63").Value = "11/11/2023 15:41"
$ws.Range("R63").Value = 2.08
$ws.Range("S63").Value = "10/11/2023 04:12"
$ws.Range("T63").Value = 1.94
$ws.Range("U63").Value = "11/11/2023 15:41"
$ws.Range("V63").Value = "https://www.betexplorer.com/football/azerbaijan/premier-league/turan-sabah-baku/8fgMRRKG/"
